{"js": "// The three \"Programa\" / \"Bibliografia\" paragraphs each contain a run whose\n// text is a single blob with no separation between numbered items (e.g.\n// \"...condutor.2) Introdu\u00e7\u00e3o...\"). The edit splits each blob into one\n// <w:t> per numbered item, joined by manual line breaks (<w:br/>) instead\n// of having them all run together.\n//\n// Word represents a manual line break inline in Range.text as \"\\u000B\" (the\n// vertical-tab code point); writing text containing \"\\u000B\" back through\n// Range.insertText reliably produces a <w:t>...</w:t><w:br/> pair in the\n// underlying OOXML, so we rebuild each paragraph's text using that\n// character as the join/separator.\n\nconst replacements = [\n  {\n    // Portuguese \"Programa\" paragraph\n    prefix: \"1) Campo Eletrost\u00e1tico e Mapeamento de Equipotenciais\",\n    items: [\n      \"1) Campo Eletrost\u00e1tico e Mapeamento de Equipotenciais: Campo de placas paralelas, Campo de cargas pontuais, Efeito de isolante e condutor.\",\n      \"2) Introdu\u00e7\u00e3o a Circuitos de Corrente Cont\u00ednua: Resistores \u00f4hmicos, Resistores n\u00e3o-ohmicos.\",\n      \"3) Resist\u00eancia e Corrente El\u00e9trica: Lei de Ohm, Modelo de Drude.\",\n      \"4) Circuitos de Corrente Cont\u00ednua: Leis de Kirchoff.\",\n      \"5) Capacitores: Associa\u00e7\u00e3o de capacitores, Carga e descarga de um capacitor.\",\n      \"6) Volt\u00edmetros, Amper\u00edmetros e Ohm\u00edmetros: Princ\u00edpio de funcionamento do Galvan\u00f4metro, Constru\u00e7\u00e3o de Volt\u00edmetros, Amper\u00edmetros e Ohm\u00edmetros.\",\n      \"7) Oscilosc\u00f3pios: Princ\u00edpio de Funcionamento do Oscilosc\u00f3pio.\",\n      \"8) Campo Magnetost\u00e1tico: Lei de Biot-Savart, Lei de Amp\u00e8re, Efeito Hall.\",\n      \"9) Lei de Indu\u00e7\u00e3o de Faraday: Indut\u00e2ncia m\u00fatua e auto-indut\u00e2ncia, Gera\u00e7\u00e3o de tens\u00e3o AC.\",\n      \"10) Circuitos RL e RC em corrente cont\u00ednua.\",\n    ],\n  },\n  {\n    // English \"Programa\" paragraph (italic run)\n    prefix: \"1) Electrostatic Field and Equipotential Mapping\",\n    items: [\n      \"1) Electrostatic Field and Equipotential Mapping: Parallel plates Field, A point charge Field, insulating effect and conductor.\",\n      \"2) Ohm\u2019s Law: ohmic resistors, resistors non-ohmic.\",\n      \"3) Resistance and Electric current: Ohm's Law, Drude model.\",\n      \"4) Direct Current Circuits: Kirchoff laws.\",\n      \"5) Capacitors: Capacitors association, load and discharge a capacitor.\",\n      \"6) Voltmeters, Ammeters and ohmmeters: Galvanometer operation principle, Voltmeters Construction, Ammeters and ohmmeters.\",\n      \"7) Oscilloscope: Oscilloscope Operation Principle.\",\n      \"8) Magnetostatic Field: Biot-Savart law, Ampere's law, Hall effect.\",\n      \"9) Faraday's Law of Induction: Mutual inductance and self-inductance, AC voltage generation.\",\n      \"10) RL and RC in DC circuits\",\n    ],\n  },\n  {\n    // \"Bibliografia\" paragraph\n    prefix: \"1. Apostilas do Laborat\u00f3rio de Ensino de F\u00edsica do IFSC/USP.\",\n    items: [\n      \"1. Apostilas do Laborat\u00f3rio de Ensino de F\u00edsica do IFSC/USP.\",\n      \"2. VUOLO, J.H. Fundamentos da Teoria de Erros, Edgard Blucher (1996).\",\n      \"3. NUSSENZVEIG, H.M. Curso de F\u00edsica B\u00e1sica. Vol. 3, Edgard Blucher (2008).\",\n      \"4. RESNICK, R.; HALLIDAY, D. Fundamentos de F\u00edsica. Vol. 3, LTC (2008).\",\n      \"5. TIPLER, P.; MOSCA, G. F\u00edsica para Cientistas e Engenheiros. Vol. 3, LTC (2008).\",\n      \"6. SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. F\u00edsica III, Vol. 3, \",\n      \"    Pearson Addison Wesley (2009).\",\n      \"7. JEWETT Jr, John W.; SERWAY, Raymond A. Princ\u00edpios de F\u00edsica. Vol. 3, Thomson Pioneira (2008).\",\n    ],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const { prefix, items } of replacements) {\n  const target = paragraphs.items.find((p) => p.text.indexOf(prefix) === 0);\n  if (!target) {\n    throw new Error(\"Paragraph not found for prefix: \" + prefix);\n  }\n  const newText = items.join(\"\\u000B\");\n  target.getRange(\"Whole\").insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The three \"Programa\" / \"Bibliografia\" paragraphs each hold a single run\n# whose text is one long blob with the numbered items run together (e.g.\n# \"...condutor.2) Introdu\u00e7\u00e3o...\"). This splits each blob into one text\n# segment per numbered item, separated by manual line breaks instead of\n# having them all flow together with no separation.\n#\n# Word (and the underlying OOXML) represents a manual line break inline in\n# Range.Text as Chr(11) (vertical tab); writing a string containing Chr(11)\n# back into a Range's .Text reliably produces a <w:t>...</w:t><w:br/> pair\n# in the saved document, so we rebuild each paragraph's text using that\n# character as the separator between items.\n\n$d = $word.ActiveDocument\n$lineBreak = [char]11\n\n$para1Items = @(\n  '1) Campo Eletrost\u00e1tico e Mapeamento de Equipotenciais: Campo de placas paralelas, Campo de cargas pontuais, Efeito de isolante e condutor.',\n  '2) Introdu\u00e7\u00e3o a Circuitos de Corrente Cont\u00ednua: Resistores \u00f4hmicos, Resistores n\u00e3o-ohmicos.',\n  '3) Resist\u00eancia e Corrente El\u00e9trica: Lei de Ohm, Modelo de Drude.',\n  '4) Circuitos de Corrente Cont\u00ednua: Leis de Kirchoff.',\n  '5) Capacitores: Associa\u00e7\u00e3o de capacitores, Carga e descarga de um capacitor.',\n  '6) Volt\u00edmetros, Amper\u00edmetros e Ohm\u00edmetros: Princ\u00edpio de funcionamento do Galvan\u00f4metro, Constru\u00e7\u00e3o de Volt\u00edmetros, Amper\u00edmetros e Ohm\u00edmetros.',\n  '7) Oscilosc\u00f3pios: Princ\u00edpio de Funcionamento do Oscilosc\u00f3pio.',\n  '8) Campo Magnetost\u00e1tico: Lei de Biot-Savart, Lei de Amp\u00e8re, Efeito Hall.',\n  '9) Lei de Indu\u00e7\u00e3o de Faraday: Indut\u00e2ncia m\u00fatua e auto-indut\u00e2ncia, Gera\u00e7\u00e3o de tens\u00e3o AC.',\n  '10) Circuitos RL e RC em corrente cont\u00ednua.'\n)\n\n$para2Items = @(\n  '1) Electrostatic Field and Equipotential Mapping: Parallel plates Field, A point charge Field, insulating effect and conductor.',\n  '2) Ohm\u2019s Law: ohmic resistors, resistors non-ohmic.',\n  '3) Resistance and Electric current: Ohm''s Law, Drude model.',\n  '4) Direct Current Circuits: Kirchoff laws.',\n  '5) Capacitors: Capacitors association, load and discharge a capacitor.',\n  '6) Voltmeters, Ammeters and ohmmeters: Galvanometer operation principle, Voltmeters Construction, Ammeters and ohmmeters.',\n  '7) Oscilloscope: Oscilloscope Operation Principle.',\n  '8) Magnetostatic Field: Biot-Savart law, Ampere''s law, Hall effect.',\n  '9) Faraday''s Law of Induction: Mutual inductance and self-inductance, AC voltage generation.',\n  '10) RL and RC in DC circuits'\n)\n\n$para3Items = @(\n  '1. Apostilas do Laborat\u00f3rio de Ensino de F\u00edsica do IFSC/USP.',\n  '2. VUOLO, J.H. Fundamentos da Teoria de Erros, Edgard Blucher (1996).',\n  '3. NUSSENZVEIG, H.M. Curso de F\u00edsica B\u00e1sica. Vol. 3, Edgard Blucher (2008).',\n  '4. RESNICK, R.; HALLIDAY, D. Fundamentos de F\u00edsica. Vol. 3, LTC (2008).',\n  '5. TIPLER, P.; MOSCA, G. F\u00edsica para Cientistas e Engenheiros. Vol. 3, LTC (2008).',\n  '6. SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. F\u00edsica III, Vol. 3, ',\n  '    Pearson Addison Wesley (2009).',\n  '7. JEWETT Jr, John W.; SERWAY, Raymond A. Princ\u00edpios de F\u00edsica. Vol. 3, Thomson Pioneira (2008).'\n)\n\nfunction Set-ParagraphByPrefix($doc, [string]$prefix, [string[]]$items, [string]$sep) {\n  $search = $doc.Content\n  $search.Find.ClearFormatting()\n  $found = $search.Find.Execute($prefix, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n  if (-not $found) {\n    throw \"Paragraph not found for prefix: $prefix\"\n  }\n  $para = $search.Paragraphs(1)\n  $newText = [string]::Join($sep, $items)\n  $para.Range.Text = $newText\n}\n\nSet-ParagraphByPrefix $d '1) Campo Eletrost\u00e1tico e Mapeamento' $para1Items $lineBreak\nSet-ParagraphByPrefix $d '1) Electrostatic Field and Equipotential Mapping' $para2Items $lineBreak\nSet-ParagraphByPrefix $d '1. Apostilas do Laborat\u00f3rio de Ensino de F\u00edsica' $para3Items $lineBreak\n"}
